$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data block for "Ajo" (garlic) records, rows 100..159, gets a new
# record inserted at the top (row 100) and every existing record shifts
# down by one row; the record that used to be last (row 159) becomes the
# new row 160. Only the date (D), volume (J), min/max/avg price (K/L/M)
# and $/Kg price (P) columns vary row to row - the rest of the columns in
# this block are constant, so we only need to shift those six columns and
# then populate the newly-appended row 160 completely (copying the
# constant columns from row 159 since row 160 didn't exist before).

$firstRow = 100
$lastRow = 159
$newLastRow = 160

# 1) Snapshot the current values of the columns that vary per record.
$dateVals = @()
$volVals = @()
$minVals = @()
$maxVals = @()
$avgVals = @()
$kgVals = @()

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $dateVals += , $ws.Cells.Item($r, 4).Value2
    $volVals += , $ws.Cells.Item($r, 10).Value2
    $minVals += , $ws.Cells.Item($r, 11).Value2
    $maxVals += , $ws.Cells.Item($r, 12).Value2
    $avgVals += , $ws.Cells.Item($r, 13).Value2
    $kgVals += , $ws.Cells.Item($r, 16).Value2
}

# 2) Create the new row 160 by copying row 159's contents cell-by-cell
#    (it is identical in every column except it now lives one row further
#    down); the six varying columns get overwritten below with the
#    shifted values. Only column D carries an explicit (date) number
#    format in this block, so that is the only one that needs copying -
#    touching NumberFormat on the other, default-styled cells would
#    needlessly fork a new style entry.
for ($c = 1; $c -le 18; $c++) {
    $src = $ws.Cells.Item(159, $c)
    $dst = $ws.Cells.Item(160, $c)
    $dst.Value = $src.Value2
}
$ws.Cells.Item(160, 4).NumberFormat = $ws.Cells.Item(159, 4).NumberFormat

# 3) Shift every existing record down by one row (160 <- old 159,
#    159 <- old 158, ... 101 <- old 100), walking from the bottom up so
#    we never clobber a value before it has been read.
for ($r = $newLastRow; $r -gt $firstRow; $r--) {
    $idx = $r - 1 - $firstRow
    $ws.Cells.Item($r, 4).Value = $dateVals[$idx]
    $ws.Cells.Item($r, 10).Value = $volVals[$idx]
    $ws.Cells.Item($r, 11).Value = $minVals[$idx]
    $ws.Cells.Item($r, 12).Value = $maxVals[$idx]
    $ws.Cells.Item($r, 13).Value = $avgVals[$idx]
    $ws.Cells.Item($r, 16).Value = $kgVals[$idx]
}

# 4) Write the brand-new record into row 100.
$ws.Cells.Item($firstRow, 4).Value = 44488
$ws.Cells.Item($firstRow, 10).Value = 250
$ws.Cells.Item($firstRow, 11).Value = 20000
$ws.Cells.Item($firstRow, 12).Value = 20000
$ws.Cells.Item($firstRow, 13).Value = 20000
$ws.Cells.Item($firstRow, 16).Value = 2000
